$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44278
$ws.Range("J2").Value = 700
$ws.Range("K2").Value = 600
$ws.Range("L2").Value = 700
$ws.Range("M2").Value = 650
$ws.Range("P2").Value = 650

# Row 3
$ws.Range("D3").Value = 44278
$ws.Range("I3").Value = "Tercera"
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 500
$ws.Range("M3").Value = 550
$ws.Range("P3").Value = 550

# Row 6
$ws.Range("D6").Value = 44224
$ws.Range("I6").Value = "Segunda"

# Row 7
$ws.Range("D7").Value = 44229
$ws.Range("J7").Value = 760
$ws.Range("K7").Value = 550
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = 575
$ws.Range("P7").Value = 575

# Row 8
$ws.Range("D8").Value = 44253
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 800
$ws.Range("L8").Value = 900
$ws.Range("M8").Value = 850
$ws.Range("P8").Value = 850

# Row 9
$ws.Range("D9").Value = 44253
$ws.Range("J9").Value = 800
$ws.Range("K9").Value = 600
$ws.Range("L9").Value = 700
$ws.Range("M9").Value = 650
$ws.Range("P9").Value = 650

# Row 10
$ws.Range("D10").Value = 44245
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 800
$ws.Range("K10").Value = 850
$ws.Range("L10").Value = 900
$ws.Range("M10").Value = 875
$ws.Range("P10").Value = 875

# Row 11
$ws.Range("D11").Value = 44245
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 750
$ws.Range("L11").Value = 800
$ws.Range("M11").Value = 775
$ws.Range("P11").Value = 775

# Row 13
$ws.Range("D13").Value = 44201
$ws.Range("J13").Value = 500

# Row 14
$ws.Range("D14").Value = 44267
$ws.Range("J14").Value = 400
$ws.Range("K14").Value = 500
$ws.Range("L14").Value = 600
$ws.Range("M14").Value = 550
$ws.Range("P14").Value = 550
